$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain plain text so numeric-looking values (e.g. "215.52",
# "1.002") are not auto-converted to floating point numbers, matching the
# original inline-string cell type used throughout the sheet.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.038.53'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = '1.645.24'
$ws.Range("E3").Value = '  -1.42%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '215.52'
$ws.Range("E5").Value = '  +2.42%  '
$ws.Range("D6").Value = '0.5220'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '0.2610'
$ws.Range("E8").Value = '  -0.29%  '
$ws.Range("D9").Value = '0.06358'
$ws.Range("E9").Value = '  +0.55%  '
$ws.Range("D10").Value = '20.83'
$ws.Range("E10").Value = '  -1.59%  '
$ws.Range("D11").Value = '0.07662'
$ws.Range("E11").Value = '  +1.59%  '
$ws.Range("D12").Value = '1.646.73'
$ws.Range("E12").Value = '  -1.44%  '
$ws.Range("D13").Value = '4.421'
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("D14").Value = '1.869.13'
$ws.Range("E14").Value = '  -1.56%  '
$ws.Range("D15").Value = '0.5543'
$ws.Range("E15").Value = '  +1.54%  '
$ws.Range("D16").Value = '0.0₅8308'
$ws.Range("E16").Value = '  +3.59%  '
$ws.Range("D17").Value = '64.88'
$ws.Range("E17").Value = '  -2.39%  '
$ws.Range("D18").Value = '26.053.51'
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").Value = '4.724'
$ws.Range("E20").Value = '  -0.49%  '
$ws.Range("D21").Value = '188.35'
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("D22").Value = '10.19'
$ws.Range("E22").Value = '  -1.09%  '
$ws.Range("D23").Value = '6.258'
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("D25").Value = '145.81'
$ws.Range("E25").Value = '  -2.62%  '
$ws.Range("D26").Value = '0.1218'
$ws.Range("E26").Value = '  -1.65%  '
$ws.Range("D27").Value = '7.411'
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("D28").Value = '15.84'
$ws.Range("E28").Value = '  +0.46%  '
$ws.Range("D29").Value = '1.401'
$ws.Range("E29").Value = '  +3.40%  '
$ws.Range("D30").Value = '0.05956'
$ws.Range("E30").Value = '  -5.36%  '
$ws.Range("D31").Value = '1.266'
$ws.Range("E31").Value = '  -1.19%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '3.407'
$ws.Range("E32").Value = '  -2.96%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '3.397'
$ws.Range("E33").Value = '  -0.72%  '
$ws.Range("D34").Value = '1.654'
$ws.Range("E34").Value = '  +0.33%  '
$ws.Range("D35").Value = '0.9972'
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("D36").Value = '2.394'
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("D37").Value = '2.754'
$ws.Range("E37").Value = '  -0.47%  '
$ws.Range("D38").Value = '0.5626'
$ws.Range("E38").Value = '  -6.45%  '
$ws.Range("D39").Value = '0.01609'
$ws.Range("E39").Value = '  -0.23%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '5.853'
$ws.Range("E40").Value = '  -3.42%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '0.8554'
$ws.Range("E41").Value = '  -0.84%  '
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").Value = '1.027.08'
$ws.Range("E43").Value = '  -7.99%  '
$ws.Range("D44").Value = '98.55'
$ws.Range("E44").Value = '  -1.97%  '
$ws.Range("D45").Value = '1.795.61'
$ws.Range("E45").Value = '  -1.44%  '
$ws.Range("E46").Value = '  +0.89%  '
$ws.Range("D47").Value = '55.71'
$ws.Range("E47").Value = '  +0.30%  '
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("D49").Value = '8.094'
$ws.Range("E49").Value = '  +0.47%  '
$ws.Range("D50").Value = '0.05151'
$ws.Range("E50").Value = '  -1.95%  '
$ws.Range("D51").Value = '0.4217'
$ws.Range("E51").Value = '  -0.50%  '
